$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "70.477.96"
Set-TextValue "E2" "  -2.91%  "
Set-TextValue "D3" "3.851.44"
Set-TextValue "E3" "  -3.18%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.13%  "
Set-TextValue "D5" "585.95"
Set-TextValue "E5" "  +0.02%  "
Set-TextValue "D6" "165.83"
Set-TextValue "E6" "  +5.46%  "
Set-TextValue "E7" "  -2.22%  "
Set-TextValue "E8" "  +0.18%  "
Set-TextValue "E9" "  -0.53%  "
Set-TextValue "E10" "  +1.76%  "
Set-TextValue "D11" "53.07"
Set-TextValue "E11" "  +0.10%  "
Set-TextValue "E12" "  -0.71%  "
Set-TextValue "D13" "11.20"
Set-TextValue "E13" "  +3.19%  "
Set-TextValue "D14" "4.460.65"
Set-TextValue "E14" "  -3.45%  "
Set-TextValue "D15" "3.886.73"
Set-TextValue "E15" "  -2.06%  "
Set-TextValue "D16" "20.60"
Set-TextValue "E16" "  +1.06%  "
Set-TextValue "D17" "13.76"
Set-TextValue "E17" "  -2.30%  "
Set-TextValue "E18" "  -6.33%  "
Set-TextValue "E19" "  -2.14%  "
Set-TextValue "D20" "70.260.17"
Set-TextValue "E20" "  -2.93%  "
Set-TextValue "D21" "434.44"
Set-TextValue "E21" "  +0.41%  "
Set-TextValue "D22" "4.67"
Set-TextValue "E22" "  -0.61%  "
Set-TextValue "D23" "93.76"
Set-TextValue "E23" "  -2.31%  "
Set-TextValue "E24" "  -4.88%  "
Set-TextValue "D25" "13.68"
Set-TextValue "E25" "  -4.49%  "
Set-TextValue "D26" "4.05"
Set-TextValue "E26" "  -8.13%  "
Set-TextValue "D27" "10.87"
Set-TextValue "E27" "  -2.04%  "
Set-TextValue "E28" "  -0.01%  "
Set-TextValue "D29" "10.21"
Set-TextValue "E29" "  -4.65%  "
Set-TextValue "D30" "34.90"
Set-TextValue "E30" "  -4.48%  "
Set-TextValue "E31" "  +1.93%  "
Set-TextValue "D32" "13.43"
Set-TextValue "E32" "  -0.99%  "
Set-TextValue "D33" "48.50"
Set-TextValue "E33" "  -0.46%  "
Set-TextValue "E34" "  -5.07%  "
Set-TextValue "D35" "69.32"
Set-TextValue "E35" "  +1.25%  "
Set-TextValue "D36" "0.0₃0968"
Set-TextValue "E36" "  +9.74%  "
Set-TextValue "D37" "614.81"
Set-TextValue "E37" "  -9.44%  "
Set-TextValue "D38" "0.416"
Set-TextValue "E38" "  -4.66%  "
Set-TextValue "D39" "1.00"
Set-TextValue "E39" "  +0.08%  "
Set-TextValue "E40" "  -0.05%  "
Set-TextValue "E41" "  -2.47%  "
Set-TextValue "E42" "  -2.98%  "
Set-TextValue "D43" "3.21"
Set-TextValue "E43" "  +25.53%  "
Set-TextValue "E44" "  -4.57%  "
Set-TextValue "D45" "9.96"
Set-TextValue "E45" "  -7.03%  "
Set-TextValue "E46" "  +0.54%  "
Set-TextValue "D47" "0.143"
Set-TextValue "E47" "  -4.05%  "
Set-TextValue "D48" "3.28"
Set-TextValue "E48" "  -2.76%  "
Set-TextValue "E49" "  -17.32%  "
Set-TextValue "D50" "2.828.71"
Set-TextValue "E50" "  +1.90%  "
Set-TextValue "D51" "0.000270"
Set-TextValue "E51" "  -0.34%  "
